$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.160.90"
$ws.Range("E2").Value = "'  -3.59%  "
$ws.Range("D3").Value = "'3.289.67"
$ws.Range("E3").Value = "'  -5.69%  "
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("D5").Value = "'541.78"
$ws.Range("E5").Value = "'  -2.29%  "
$ws.Range("D6").Value = "'169.70"
$ws.Range("E6").Value = "'  -4.95%  "
$ws.Range("D7").Value = "'0.609"
$ws.Range("E7").Value = "'  -4.56%  "
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("D9").Value = "'3.276.96"
$ws.Range("E9").Value = "'  -5.78%  "
$ws.Range("E10").Value = "'  -4.44%  "
$ws.Range("E11").Value = "'  -1.02%  "
$ws.Range("D12").Value = "'52.38"
$ws.Range("E12").Value = "'  -2.41%  "
$ws.Range("D13").Value = "'0.0000262"
$ws.Range("E13").Value = "'  -2.99%  "
$ws.Range("D14").Value = "'8.77"
$ws.Range("E14").Value = "'  -5.16%  "
$ws.Range("D15").Value = "'3.827.02"
$ws.Range("E15").Value = "'  -5.38%  "
$ws.Range("D16").Value = "'17.81"
$ws.Range("E16").Value = "'  -4.68%  "
$ws.Range("B17").Value = "'TRON"
$ws.Range("C17").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.116"
$ws.Range("E17").Value = "'  -4.35%  "
$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.292.53"
$ws.Range("E18").Value = "'  -5.58%  "
$ws.Range("D19").Value = "'11.54"
$ws.Range("E19").Value = "'  -4.72%  "
$ws.Range("D20").Value = "'63.036.97"
$ws.Range("E20").Value = "'  -3.77%  "
$ws.Range("D21").Value = "'0.963"
$ws.Range("E21").Value = "'  -2.67%  "
$ws.Range("D22").Value = "'411.72"
$ws.Range("E22").Value = "'  -1.38%  "
$ws.Range("B23").Value = "'Toncoin"
$ws.Range("C23").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'4.36"
$ws.Range("E23").Value = "'  +5.78%  "
$ws.Range("B24").Value = "'PancakeSwap"
$ws.Range("C24").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'4.00"
$ws.Range("E24").Value = "'  -1.00%  "
$ws.Range("D25").Value = "'13.56"
$ws.Range("E25").Value = "'  +5.36%  "
$ws.Range("D26").Value = "'82.29"
$ws.Range("E26").Value = "'  -4.55%  "
$ws.Range("D27").Value = "'10.43"
$ws.Range("E27").Value = "'  -3.36%  "
$ws.Range("E28").Value = "'  -5.43%  "
$ws.Range("D29").Value = "'8.53"
$ws.Range("E29").Value = "'  -5.69%  "
$ws.Range("D30").Value = "'28.73"
$ws.Range("E30").Value = "'  -4.99%  "
$ws.Range("D31").Value = "'6.29"
$ws.Range("E31").Value = "'  -3.28%  "
$ws.Range("D32").Value = "'11.23"
$ws.Range("E32").Value = "'  -4.28%  "
$ws.Range("D33").Value = "'564.51"
$ws.Range("E33").Value = "'  -7.66%  "
$ws.Range("E34").Value = "'  -4.12%  "
$ws.Range("D35").Value = "'57.48"
$ws.Range("E35").Value = "'  -3.33%  "
$ws.Range("E36").Value = "'  +0.03%  "
$ws.Range("E37").Value = "'  -1.35%  "
$ws.Range("D38").Value = "'34.76"
$ws.Range("E38").Value = "'  -7.15%  "
$ws.Range("E39").Value = "'  +4.49%  "
$ws.Range("D40").Value = "'0.0₃0729"
$ws.Range("E40").Value = "'  -7.18%  "
$ws.Range("E41").Value = "'  -5.08%  "
$ws.Range("D42").Value = "'3.096.41"
$ws.Range("E42").Value = "'  -8.73%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "'  +0.08%  "
$ws.Range("D44").Value = "'2.73"
$ws.Range("E44").Value = "'  -3.01%  "
$ws.Range("E45").Value = "'  -0.87%  "
$ws.Range("D46").Value = "'0.0396"
$ws.Range("E46").Value = "'  -4.31%  "
$ws.Range("E47").Value = "'  -5.86%  "
$ws.Range("E48").Value = "'  -4.15%  "
$ws.Range("E49").Value = "'  -4.31%  "
$ws.Range("D50").Value = "'131.73"
$ws.Range("E50").Value = "'  -4.29%  "
$ws.Range("E51").Value = "'  -6.51%  "
